$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9803.076596270917
$ws.Range("D2").Value = -1.623116802729121
$ws.Range("E2").Value = -13.11013451013214
$ws.Range("F2").Value = -0.2215765152623252
$ws.Range("G2").Value = 7.042153032673656
$ws.Range("H2").Value = -0.1408537640523785
$ws.Range("I2").Value = 4.014994584563411
$ws.Range("J2").Value = -0.06827593745902204
$ws.Range("K2").Value = 1.078454503891648
$ws.Range("L2").Value = 0.01485274608893764
$ws.Range("M2").Value = -0.1107468848678981
$ws.Range("N2").Value = -0.003210158682024669
$ws.Range("O2").Value = 0.01576594610035043
$ws.Range("P2").Value = -0.03372707259556925
$ws.Range("Q2").Value = -0.004786091951473356
$ws.Range("R2").Value = -85.85462730574083
$ws.Range("S2").Value = -0.000000000000004056042687253865
$ws.Range("T2").Value = 0.01083374460524887
$ws.Range("U2").Value = -0.0000000000000004780929021570105
$ws.Range("V2").Value = 0.00001168890917761487
$ws.Range("W2").Value = 0.000000000000007752256655194196
$ws.Range("X2").Value = 0.000000003021160356279676
$ws.Range("Y2").Value = -0.000000000000002296035900163013
$ws.Range("Z2").Value = 0.000000000003412687757586493
$ws.Range("AA2").Value = -0.00000000000000299450830292504
$ws.Range("AB2").Value = 0.000000000000007343463381509304
$ws.Range("AC2").Value = 0.00000000000001122228364368286
$ws.Range("AD2").Value = 0.000000000000001272973386458379
$ws.Range("AE2").Value = 0.000000000000009585799867875783
$ws.Range("AF2").Value = 0.00000000000003860851780485226
$ws.Range("C3").Value = 9898.075730671611
$ws.Range("D3").Value = 0.3450018307690805
$ws.Range("E3").Value = -19.77726182624921
$ws.Range("F3").Value = -0.2149133672519221
$ws.Range("G3").Value = 1.454761990727579
$ws.Range("H3").Value = -0.2114221760266576
$ws.Range("I3").Value = 4.488689910101072
$ws.Range("J3").Value = -0.09077645461035043
$ws.Range("K3").Value = 1.150603747240015
$ws.Range("L3").Value = 0.0008290466532624503
$ws.Range("M3").Value = -0.07163832561740857
$ws.Range("N3").Value = -0.01107789830231555
$ws.Range("O3").Value = -0.2692574943686147
$ws.Range("P3").Value = -0.02001654987988458
$ws.Range("Q3").Value = -0.002832155302489597
$ws.Range("R3").Value = 4250.736354136774
$ws.Range("S3").Value = 0.0000000000001899317054832392
$ws.Range("T3").Value = -4.904374331053778
$ws.Range("U3").Value = -0.0000000000001781444734237937
$ws.Range("V3").Value = -0.6781251582048761
$ws.Range("W3").Value = -0.0000000000002989455804501879
$ws.Range("X3").Value = 0.944149198220479
$ws.Range("Y3").Value = 0.0000000000003662296367272287
$ws.Range("Z3").Value = 0.1810032173126672
$ws.Range("AA3").Value = -0.0000000000008257709569998079
$ws.Range("AB3").Value = -0.04628955808208494
$ws.Range("AC3").Value = 0.0000000000005598935723348323
$ws.Range("AD3").Value = -0.03146475895843957
$ws.Range("AE3").Value = 0.0000000000002519794515074252
$ws.Range("AF3").Value = 0.0003016796947173021
$ws.Range("C4").Value = 10000
$ws.Range("D4").Value = 1.070914692675539
$ws.Range("E4").Value = -19.85378874922338
$ws.Range("F4").Value = -0.2137313394331523
$ws.Range("G4").Value = 1.140489881099531
$ws.Range("H4").Value = -0.2500210138594487
$ws.Range("I4").Value = 4.633856183512828
$ws.Range("J4").Value = -0.08980684038600628
$ws.Range("K4").Value = 1.15741839350241
$ws.Range("L4").Value = 0.006762133757860266
$ws.Range("M4").Value = -0.07729436521041172
$ws.Range("N4").Value = 0.05606664521192392
$ws.Range("O4").Value = -0.1806724710233127
$ws.Range("P4").Value = -0.01947884215303555
$ws.Range("Q4").Value = 0.010495978630842
$ws.Range("R4").Value = 10053.83512543139
$ws.Range("S4").Value = -0.000000000001364798278083048
$ws.Range("T4").Value = -16.11576858065114
$ws.Range("U4").Value = 0.0000000000004060017147739935
$ws.Range("V4").Value = 1.986617295214211
$ws.Range("W4").Value = 0.0000000000005973232021978054
$ws.Range("X4").Value = 5.757784810706341
$ws.Range("Y4").Value = -0.0000000000006420715676951326
$ws.Range("Z4").Value = 1.155825181643935
$ws.Range("AA4").Value = -0.00000000000003831427268733069
$ws.Range("AB4").Value = -0.0403430232133415
$ws.Range("AC4").Value = 0.000000000001029451608597849
$ws.Range("AD4").Value = -0.1046127229087321
$ws.Range("AE4").Value = -0.0000000000007495420215292956
$ws.Range("AF4").Value = 0.002594381419169133
$ws.Range("C5").Value = 9905.314281694866
$ws.Range("D5").Value = 1.487355479645196
$ws.Range("E5").Value = -19.85786241493973
$ws.Range("F5").Value = -0.1166025697568185
$ws.Range("G5").Value = 1.362708154051409
$ws.Range("H5").Value = -0.1707759698383113
$ws.Range("I5").Value = 4.523356416628299
$ws.Range("J5").Value = -0.09870713713772904
$ws.Range("K5").Value = 1.164895871930481
$ws.Range("L5").Value = 0.006898238153751401
$ws.Range("M5").Value = -0.077803920052874
$ws.Range("N5").Value = 0.01473090747404357
$ws.Range("O5").Value = -0.1905472933990691
$ws.Range("P5").Value = -0.04494949079791635
$ws.Range("Q5").Value = 0.02950162784422101
$ws.Range("R5").Value = 4260.813524873376
$ws.Range("S5").Value = -0.0000000000005820229535632692
$ws.Range("T5").Value = -1.392464301102987
$ws.Range("U5").Value = -0.00000000000006738650305857813
$ws.Range("V5").Value = -0.05916572131225438
$ws.Range("W5").Value = 0.0000000000008780859543495461
$ws.Range("X5").Value = 1.087766970577664
$ws.Range("Y5").Value = -0.0000000000004187554964396916
$ws.Range("Z5").Value = 0.2077226363640378
$ws.Range("AA5").Value = 0.0000000000007790987287185776
$ws.Range("AB5").Value = -0.04144024203183552
$ws.Range("AC5").Value = 0.0000000000005381684157529254
$ws.Range("AD5").Value = -0.03084241377600528
$ws.Range("AE5").Value = 0.0000000000005957537801052465
$ws.Range("AF5").Value = 0.0003894497799132393
$ws.Range("C6").Value = 10038.68010434441
$ws.Range("D6").Value = 1.318081070584686
$ws.Range("E6").Value = -19.94638767334052
$ws.Range("F6").Value = -0.03584768603149048
$ws.Range("G6").Value = 1.569503487914278
$ws.Range("H6").Value = -0.2156058750838234
$ws.Range("I6").Value = 4.51255664097599
$ws.Range("J6").Value = -0.114848492319917
$ws.Range("K6").Value = 1.17402908713419
$ws.Range("L6").Value = 0.01254728869608168
$ws.Range("M6").Value = -0.08931624690410626
$ws.Range("N6").Value = -0.01620071842649914
$ws.Range("O6").Value = -0.2633180167585202
$ws.Range("P6").Value = 0.004169609869314606
$ws.Range("Q6").Value = 0.00693944031204616
$ws.Range("R6").Value = -83.50427399563998
$ws.Range("S6").Value = -0.00000000000001414880243919718
$ws.Range("T6").Value = 0.01908398409562601
$ws.Range("U6").Value = -0.000000000000001509290685128469
$ws.Range("V6").Value = 0.0001045148195935405
$ws.Range("W6").Value = 0.00000000000001073174698991082
$ws.Range("X6").Value = 0.000001432596819505779
$ws.Range("Y6").Value = -0.000000000000002714130288343722
$ws.Range("Z6").Value = 0.00000002544464126453821
$ws.Range("AA6").Value = -0.000000000000007277983347676765
$ws.Range("AB6").Value = 0.0000000004951644862091057
$ws.Range("AC6").Value = -0.00000000000001221386055043953
$ws.Range("AD6").Value = 0.00000000001024059996922259
$ws.Range("AE6").Value = 0.000000000000003750290672352939
$ws.Range("AF6").Value = 0.000000000000216149271736694
